# New weekly price record for "Vega Modelo de Temuco - Ciboulette" was
# inserted as a new data row (362), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 362, shifting rows 362:444
# down to 363:445 (dimension grows from A1:R444 to A1:R445).
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(362, 1).Value = 10
$ws.Cells.Item(362, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(362, 3).Value = "La Araucanía"
$ws.Cells.Item(362, 4).Value = 45244
$ws.Cells.Item(362, 5).Value = 9
$ws.Cells.Item(362, 6).Value = 100112039
$ws.Cells.Item(362, 7).Value = "Ciboulette"
$ws.Cells.Item(362, 8).Value = "Sin especificar"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 65
$ws.Cells.Item(362, 11).Value = 7000
$ws.Cells.Item(362, 12).Value = 7000
$ws.Cells.Item(362, 13).Value = 7000
$ws.Cells.Item(362, 14).Value = "`$/docena de atados"
$ws.Cells.Item(362, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(362, 16).Value = 2333
$ws.Cells.Item(362, 17).Value = 3
$ws.Cells.Item(362, 18).Value = "Hortaliza"

# Match the date-number-format already used by the other rows in column D.
$ws.Cells.Item(362, 4).NumberFormat = $ws.Cells.Item(363, 4).NumberFormat
